# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to rows across the 8 craft sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Cells.Item(132, 8).Value = 5032.852
$ws.Cells.Item(132, 9).Value = 4722.1816
$ws.Cells.Item(132, 10).Value = 6399.8
$ws.Cells.Item(132, 11).Value = 14166.5448
$ws.Cells.Item(132, 12).Value = 19199.4
$ws.Cells.Item(132, 13).Value = -11636.5448
$ws.Cells.Item(132, 14).Value = -24259.4

# Row 137
$ws.Cells.Item(137, 8).Value = 1815.875
$ws.Cells.Item(137, 9).Value = 1585.8572
$ws.Cells.Item(137, 10).Value = 1994.7778
$ws.Cells.Item(137, 11).Value = 4757.571599999999
$ws.Cells.Item(137, 12).Value = 5984.3334
$ws.Cells.Item(137, 13).Value = -2207.571599999999
$ws.Cells.Item(137, 14).Value = -11084.3334

# Row 138
$ws.Cells.Item(138, 8).Value = 3120.0256
$ws.Cells.Item(138, 10).Value = 3061.1785
$ws.Cells.Item(138, 12).Value = 9183.5355
$ws.Cells.Item(138, 14).Value = -19463.5355

$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Cells.Item(28, 8).Value = 22485.375
$ws.Cells.Item(28, 9).Value = 6820.3335
$ws.Cells.Item(28, 10).Value = 31884.4
$ws.Cells.Item(28, 11).Value = 6820.3335
$ws.Cells.Item(28, 12).Value = 31884.4
$ws.Cells.Item(28, 13).Value = -6628.3335
$ws.Cells.Item(28, 14).Value = -32268.4

# Row 45
$ws.Cells.Item(45, 8).Value = 3226.2354
$ws.Cells.Item(45, 9).Value = 2263
$ws.Cells.Item(45, 10).Value = 4082.4443
$ws.Cells.Item(45, 11).Value = 2263
$ws.Cells.Item(45, 12).Value = 4082.4443
$ws.Cells.Item(45, 13).Value = -1886
$ws.Cells.Item(45, 14).Value = -4836.4443

# Row 74
$ws.Cells.Item(74, 8).Value = 717.4054
$ws.Cells.Item(74, 9).Value = 569.8125
$ws.Cells.Item(74, 10).Value = 829.8570999999999
$ws.Cells.Item(74, 11).Value = 569.8125
$ws.Cells.Item(74, 12).Value = 829.8570999999999
$ws.Cells.Item(74, 13).Value = 304.1875
$ws.Cells.Item(74, 14).Value = -2577.8571

# Row 77
$ws.Cells.Item(77, 8).Value = 717.4054
$ws.Cells.Item(77, 9).Value = 569.8125
$ws.Cells.Item(77, 10).Value = 829.8570999999999
$ws.Cells.Item(77, 11).Value = 2849.0625
$ws.Cells.Item(77, 12).Value = 4149.2855
$ws.Cells.Item(77, 13).Value = 1518.9375
$ws.Cells.Item(77, 14).Value = -12885.2855

# Row 99
$ws.Cells.Item(99, 8).Value = 22485.375
$ws.Cells.Item(99, 9).Value = 6820.3335
$ws.Cells.Item(99, 10).Value = 31884.4
$ws.Cells.Item(99, 11).Value = 6820.3335
$ws.Cells.Item(99, 12).Value = 31884.4
$ws.Cells.Item(99, 13).Value = -3825.3335
$ws.Cells.Item(99, 14).Value = -37874.4

# Row 132
$ws.Cells.Item(132, 8).Value = 3562.628
$ws.Cells.Item(132, 9).Value = 3123.4827
$ws.Cells.Item(132, 10).Value = 4472.2856
$ws.Cells.Item(132, 11).Value = 9370.4481
$ws.Cells.Item(132, 12).Value = 13416.8568
$ws.Cells.Item(132, 13).Value = -6840.4481
$ws.Cells.Item(132, 14).Value = -18476.8568

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1722.5862
$ws.Cells.Item(20, 9).Value = 1934.125
$ws.Cells.Item(20, 10).Value = 1462.2307
$ws.Cells.Item(20, 11).Value = 1934.125
$ws.Cells.Item(20, 12).Value = 1462.2307
$ws.Cells.Item(20, 13).Value = -1687.125
$ws.Cells.Item(20, 14).Value = -1956.2307

# Row 64
$ws.Cells.Item(64, 8).Value = 1035.1428
$ws.Cells.Item(64, 9).Value = 923
$ws.Cells.Item(64, 10).Value = 1119.25
$ws.Cells.Item(64, 11).Value = 923
$ws.Cells.Item(64, 12).Value = 1119.25
$ws.Cells.Item(64, 13).Value = -698
$ws.Cells.Item(64, 14).Value = -1569.25

# Row 67
$ws.Cells.Item(67, 8).Value = 1035.1428
$ws.Cells.Item(67, 9).Value = 923
$ws.Cells.Item(67, 10).Value = 1119.25
$ws.Cells.Item(67, 11).Value = 923
$ws.Cells.Item(67, 12).Value = 1119.25
$ws.Cells.Item(67, 13).Value = -143
$ws.Cells.Item(67, 14).Value = -2679.25

# Row 134
$ws.Cells.Item(134, 8).Value = 2277.537
$ws.Cells.Item(134, 9).Value = 2112.366
$ws.Cells.Item(134, 10).Value = 2798.4614
$ws.Cells.Item(134, 11).Value = 6337.098
$ws.Cells.Item(134, 12).Value = 8395.3842
$ws.Cells.Item(134, 13).Value = -3802.098
$ws.Cells.Item(134, 14).Value = -13465.3842

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 3450
$ws.Cells.Item(62, 9).Value = 3450
$ws.Cells.Item(62, 11).Value = 3450
$ws.Cells.Item(62, 13).Value = -2826

# Row 65
$ws.Cells.Item(65, 8).Value = 3450
$ws.Cells.Item(65, 9).Value = 3450
$ws.Cells.Item(65, 11).Value = 17250
$ws.Cells.Item(65, 13).Value = -14130

# Row 105
$ws.Cells.Item(105, 8).Value = 325.77777
$ws.Cells.Item(105, 9).Value = 325.77777
$ws.Cells.Item(105, 11).Value = 325.77777
$ws.Cells.Item(105, 13).Value = 1421.22223

# Row 135
$ws.Cells.Item(135, 8).Value = 51452
$ws.Cells.Item(135, 10).Value = 51452
$ws.Cells.Item(135, 12).Value = 51452
$ws.Cells.Item(135, 14).Value = -61592

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 732.4020400000001
$ws.Cells.Item(68, 9).Value = 558.96295
$ws.Cells.Item(68, 10).Value = 799.3
$ws.Cells.Item(68, 11).Value = 1676.88885
$ws.Cells.Item(68, 12).Value = 2397.9
$ws.Cells.Item(68, 13).Value = -865.8888499999998
$ws.Cells.Item(68, 14).Value = -4019.9

# Row 71
$ws.Cells.Item(71, 8).Value = 732.4020400000001
$ws.Cells.Item(71, 9).Value = 558.96295
$ws.Cells.Item(71, 10).Value = 799.3
$ws.Cells.Item(71, 11).Value = 5030.66655
$ws.Cells.Item(71, 12).Value = 7193.7
$ws.Cells.Item(71, 13).Value = -974.6665499999999
$ws.Cells.Item(71, 14).Value = -15305.7

# Row 107
$ws.Cells.Item(107, 8).Value = 1377.1515
$ws.Cells.Item(107, 9).Value = 255.67647
$ws.Cells.Item(107, 11).Value = 767.02941
$ws.Cells.Item(107, 13).Value = 1152.97059

# Row 113
$ws.Cells.Item(113, 8).Value = 596.9796
$ws.Cells.Item(113, 9).Value = 387.15625
$ws.Cells.Item(113, 10).Value = 991.94116
$ws.Cells.Item(113, 11).Value = 1161.46875
$ws.Cells.Item(113, 12).Value = 2975.82348
$ws.Cells.Item(113, 13).Value = 1008.53125
$ws.Cells.Item(113, 14).Value = -7315.82348

# Row 137
$ws.Cells.Item(137, 8).Value = 11304.143
$ws.Cells.Item(137, 9).Value = 16232.25
$ws.Cells.Item(137, 10).Value = 4733.3335
$ws.Cells.Item(137, 11).Value = 48696.75
$ws.Cells.Item(137, 12).Value = 14200.0005
$ws.Cells.Item(137, 13).Value = -43596.75
$ws.Cells.Item(137, 14).Value = -24400.0005

# Row 140
$ws.Cells.Item(140, 8).Value = 2281.0527
$ws.Cells.Item(140, 9).Value = 1720
$ws.Cells.Item(140, 11).Value = 5160
$ws.Cells.Item(140, 13).Value = 20

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Cells.Item(93, 8).Value = 116666.664
$ws.Cells.Item(93, 10).Value = 116666.664
$ws.Cells.Item(93, 12).Value = 116666.664
$ws.Cells.Item(93, 14).Value = -120410.664

# Row 108
$ws.Cells.Item(108, 8).Value = 80342
$ws.Cells.Item(108, 10).Value = 80342
$ws.Cells.Item(108, 12).Value = 80342
$ws.Cells.Item(108, 14).Value = -88022

# Row 132
$ws.Cells.Item(132, 8).Value = 2473.818
$ws.Cells.Item(132, 9).Value = 2235.4075
$ws.Cells.Item(132, 10).Value = 2852.4707
$ws.Cells.Item(132, 11).Value = 6706.2225
$ws.Cells.Item(132, 12).Value = 8557.4121
$ws.Cells.Item(132, 13).Value = -4176.2225
$ws.Cells.Item(132, 14).Value = -13617.4121

$ws = $wb.Worksheets.Item("LTW")
# Row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()

# Row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()

# Row 136
$ws.Cells.Item(136, 8).Value = 1415.0454
$ws.Cells.Item(136, 9).Value = 1400.9412
$ws.Cells.Item(136, 10).Value = 1463
$ws.Cells.Item(136, 11).Value = 4202.8236
$ws.Cells.Item(136, 12).Value = 4389
$ws.Cells.Item(136, 13).Value = -1652.8236
$ws.Cells.Item(136, 14).Value = -9489

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Cells.Item(136, 8).Value = 2726.3262
$ws.Cells.Item(136, 9).Value = 2625.3928
$ws.Cells.Item(136, 11).Value = 7876.178400000001
$ws.Cells.Item(136, 13).Value = -5326.178400000001
